$wb = $excel.ActiveWorkbook

# Duplicate the original sheet (the copy inherits the worksheet's XML
# namespaces/ignorable markup) so the new sheet picks up sheetId 2, then
# drop the original Sheet1 so the duplicate becomes the sole/active sheet -
# this mirrors the target workbook.xml
# (<sheet name="ValidLogin" sheetId="2" r:id="rId1"/>).
$old = $wb.ActiveSheet
$oldName = $old.Name
$old.Copy($null, $old)
[void]$wb.Worksheets.Item($oldName).Delete()
$ws = $wb.Worksheets.Item($oldName + " (2)")
$ws.Name = "ValidLogin"
[void]$ws.Activate()

# Set up the header row and data row for the login test data
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# Update view settings to match the target state
$excel.ActiveWindow.Zoom = 175
[void]$ws.Range("B3").Select()
